# Fix matched counts in the multi-summary stream report layout.
# On both the "summary" and "revsummary" sheets, the B/H value columns in
# the "(Source1)/(Source2) - Matching/Non-matching Rows" rows (for both the
# "Inner vs Outer"-style row block 3-13 and the "Left vs Right"-style row
# block 17-27) had their counts swapped/incorrect. Correct them here.

$wb = $excel.ActiveWorkbook

$sheetNames = @("summary", "revsummary")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Rows 3-13 block (columns B and H)
    $ws.Range("B7").Value = 39.0
    $ws.Range("B9").Value = 960.0
    $ws.Range("B12").Value = 960.0

    $ws.Range("H7").Value = 0.0
    $ws.Range("H9").Value = 960.0
    $ws.Range("H12").Value = 960.0

    # Rows 17-27 block (columns B and H)
    $ws.Range("B21").Value = 7.0
    $ws.Range("B23").Value = 960.0
    $ws.Range("B26").Value = 960.0

    $ws.Range("H21").Value = 32.0
    $ws.Range("H23").Value = 960.0
    $ws.Range("H26").Value = 960.0
}
